$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 17:22"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 404056
$ws.Range("C4").Value = 3721
$ws.Range("D4").Value = 21815
$ws.Range("E4").Value = 369253
$ws.Range("F4").Value = 9200
$ws.Range("G4").Value = 147
$ws.Range("H4").Value = 12988

# Row 16: Canada - refreshed totals
$ws.Range("B16").Value = 18447
$ws.Range("C16").Value = 550
$ws.Range("E16").Value = 14016

# Row 19: Austria - refreshed totals
$ws.Range("B19").Value = 12901
$ws.Range("C19").Value = 262
$ws.Range("E19").Value = 8116

# Row 50: Grecia - refreshed totals
$ws.Range("B50").Value = 1884
$ws.Range("C50").Value = 52
$ws.Range("E50").Value = 1532
$ws.Range("F50").Value = 84
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 83

# Row 64: now Moldavia (jumped above Eslovenia/Bielorrusia with new totals)
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 1174
$ws.Range("C64").Value = 118
$ws.Range("D64").Value = 40
$ws.Range("E64").Value = 1107
$ws.Range("F64").Value = 80
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 27

# Row 65: now Eslovenia (shifted down one row, values unchanged)
$ws.Range("A65").Value = "Eslovenia"
$ws.Range("B65").Value = 1091
$ws.Range("C65").Value = 32
$ws.Range("D65").Value = 120
$ws.Range("E65").Value = 931
$ws.Range("F65").Value = 35
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 40

# Row 66: now Bielorrusia (shifted down one row, values unchanged)
$ws.Range("A66").Value = "Bielorrusia"
$ws.Range("B66").Value = 1066
$ws.Range("C66").Value = 205
$ws.Range("D66").Value = 77
$ws.Range("E66").Value = 976
$ws.Range("F66").Value = 33
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 13

# Row 86: Republica de Chipre - refreshed totals
$ws.Range("B86").Value = 526
$ws.Range("C86").Value = 32
$ws.Range("D86").Value = 52
$ws.Range("E86").Value = 465

# Row 98: now Ghana (jumped above Honduras with new totals)
$ws.Range("A98").Value = "Ghana"
$ws.Range("B98").Value = 313
$ws.Range("C98").Value = 26
$ws.Range("D98").Value = 31
$ws.Range("E98").Value = 277
$ws.Range("F98").Value = 2
$ws.Range("H98").Value = 5

# Row 99: now Honduras (shifted down one row, values unchanged)
$ws.Range("A99").Value = "Honduras"
$ws.Range("B99").Value = 312
$ws.Range("C99").Value = 7
$ws.Range("D99").Value = 6
$ws.Range("E99").Value = 284
$ws.Range("F99").Value = 10
$ws.Range("H99").Value = 22

# Row 100: now Malta (shifted down one row, values unchanged)
$ws.Range("A100").Value = "Malta"
$ws.Range("B100").Value = 299
$ws.Range("C100").Value = 6
$ws.Range("D100").Value = 5
$ws.Range("E100").Value = 294
$ws.Range("F100").Value = 4
$ws.Range("H100").Value = 0

# Row 119: Isla de Man - refreshed totals
$ws.Range("D119").Value = 82
$ws.Range("E119").Value = 75

# Row 132: Guatemala - refreshed totals
$ws.Range("B132").Value = 87
$ws.Range("C132").Value = 10
$ws.Range("E132").Value = 67
